# Updated cryptos list on Sat Sep  7 22:32:44 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to stay a text string even
# when it looks like a number (e.g. "127.00" would otherwise be stored as
# the number 127 and lose its trailing zero). We flip the cell to text
# format just long enough to assign the value, then restore the "Normal"
# cell style so no stray number-format is left behind on the cell.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "53.941.70"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.249.92"
$ws.Range("E3").Value = "  +2.39%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "492.84"
$ws.Range("E5").Value = "  +1.61%  "

# Row 6 - Solana
Set-TextValue "D6" "127.00"
$ws.Range("E6").Value = "  +1.69%  "

# Row 7 - USDC
Set-TextValue "D7" "0.997"
$ws.Range("E7").Value = "  +0.16%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.85%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.0951"
$ws.Range("E9").Value = "  +3.59%  "

# Row 10 - TRON
Set-TextValue "D10" "0.152"
$ws.Range("E10").Value = "  +2.54%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +3.50%  "

# Row 12 - Toncoin
Set-TextValue "D12" "4.63"
$ws.Range("E12").Value = "  +0.24%  "

# Row 13 - Wrapped liquid staked Ether 2.0
Set-TextValue "D13" "2.651.61"
$ws.Range("E13").Value = "  +2.38%  "

# Row 14 - Avalanche
Set-TextValue "D14" "21.69"
$ws.Range("E14").Value = "  +2.60%  "

# Row 15 - Wrapped BTC
Set-TextValue "D15" "53.881.25"
$ws.Range("E15").Value = "  +0.73%  "

# Row 16 - Shiba Inu
$ws.Range("E16").Value = "  +0.62%  "

# Row 17 - Wrapped Ether
Set-TextValue "D17" "2.260.92"
$ws.Range("E17").Value = "  +2.98%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  +4.29%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +3.19%  "

# Row 20 & 21 - Uniswap and BitcoinCash swapped rank positions
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "299.56"
$ws.Range("E20").Value = "  +1.77%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D21" "6.41"
$ws.Range("E21").Value = "  +5.81%  "

# Row 22 - Dai
Set-TextValue "D22" "0.999"
$ws.Range("E22").Value = "  +0.19%  "

# Row 23 - LEO
Set-TextValue "D23" "5.40"
$ws.Range("E23").Value = "  -1.88%  "

# Row 24 - Litecoin
$ws.Range("E24").Value = "  -0.96%  "

# Row 25 - Binance-Peg BSC-USD
$ws.Range("E25").Value = "  +2.06%  "

# Row 26 - Polygon
$ws.Range("E26").Value = "  +1.13%  "

# Row 27 - Wrapped eETH
Set-TextValue "D27" "2.358.08"
$ws.Range("E27").Value = "  +2.63%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +1.63%  "

# Row 29 - Internet Computer (DFINITY)
Set-TextValue "D29" "7.01"
$ws.Range("E29").Value = "  +0.32%  "

# Row 30 - Monero
Set-TextValue "D30" "165.40"
$ws.Range("E30").Value = "  +0.44%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.59"
$ws.Range("E31").Value = "  +1.17%  "

# Row 32 - PEPE
Set-TextValue "D32" "0.0₃0677"
$ws.Range("E32").Value = "  +2.20%  "

# Row 34 - Aptos
Set-TextValue "D34" "5.81"
$ws.Range("E34").Value = "  +2.60%  "

# Row 35 - First Digital USD
Set-TextValue "D35" "0.996"
$ws.Range("E35").Value = "  +0.25%  "

# Row 36 - Fetch.AI
$ws.Range("E36").Value = "  -0.41%  "

# Row 37 - Ethereum Classic
Set-TextValue "D37" "17.60"
$ws.Range("E37").Value = "  +1.71%  "

# Row 38 - Sui Network
Set-TextValue "D38" "0.883"
$ws.Range("E38").Value = "  +7.21%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +2.58%  "

# Row 40 - NEAR Protocol
Set-TextValue "D40" "3.65"
$ws.Range("E40").Value = "  +3.20%  "

# Row 41 - OKB
Set-TextValue "D41" "35.69"
$ws.Range("E41").Value = "  -0.09%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  +2.21%  "

# Row 43 - Polygon Ecosystem Token
Set-TextValue "D43" "0.370"
$ws.Range("E43").Value = "  +1.00%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  +2.30%  "

# Row 45 - Render Token
$ws.Range("E45").Value = "  +3.17%  "

# Row 46 - Aave
Set-TextValue "D46" "124.73"
$ws.Range("E46").Value = "  +0.18%  "

# Row 47 - Stellar
Set-TextValue "D47" "0.0885"
$ws.Range("E47").Value = "  +1.30%  "

# Row 48 - Mantle
$ws.Range("E48").Value = "  +1.40%  "

# Row 49 - Bittensor
Set-TextValue "D49" "235.01"
$ws.Range("E49").Value = "  +1.16%  "

# Row 50 - Hedera
$ws.Range("E50").Value = "  +2.71%  "

# Row 51 - VeChain
Set-TextValue "D51" "0.0200"
$ws.Range("E51").Value = "  +0.27%  "
